$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.936.61"
$ws.Range("E2").Value = "  +2.76%  "

$ws.Range("D3").Value = "3.087.43"
$ws.Range("E3").Value = "  +4.97%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.21%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.76"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.50%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").Value = "3.083.86"
$ws.Range("E8").Value = "  +4.96%  "

$ws.Range("E9").Value = "  +1.31%  "

$ws.Range("E10").Value = "  -1.27%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.155"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.41%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.482"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.51%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000250"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.51%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.47"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.41%  "

$ws.Range("B16").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C16").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").Value = "3.596.40"
$ws.Range("E16").Value = "  +4.83%  "

$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "66.894.53"
$ws.Range("E17").Value = "  +2.68%  "

$ws.Range("E18").Value = "  +3.74%  "

$ws.Range("D19").Value = "3.086.48"
$ws.Range("E19").Value = "  +5.03%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +8.85%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "466.44"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.79%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.715"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.19%  "

$ws.Range("E23").Value = "  +4.22%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.35"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.48%  "

$ws.Range("E25").Value = "  +7.33%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.96"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +7.37%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.15"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.39%  "

$ws.Range("E28").Value = "  +0.01%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.99"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.21%  "

$ws.Range("E30").Value = "  +0.76%  "

$ws.Range("E31").Value = "  +3.66%  "

$ws.Range("E32").Value = "  +1.15%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "28.21"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.03%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.115"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.72%  "

$ws.Range("E35").Value = "  -0.04%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.37%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.89"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.92%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.11"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.93%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "47.08"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.42%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.320"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.71%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "50.18"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.23%  "

$ws.Range("E42").Value = "  +2.05%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.68"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.68%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.83"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.08%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0361"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.94%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "384.27"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.28%  "

$ws.Range("D47").Value = "2.768.10"
$ws.Range("E47").Value = "  +2.35%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "134.93"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.47%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.79"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.93%  "

$ws.Range("E51").Value = "  +2.41%  "
